$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 250, pushing existing rows 250-330 down to 251-331
$ws.Rows.Item(250).Insert()

# Fill the new row 250 with its values
$ws.Cells.Item(250, 1).Value = 6
$ws.Cells.Item(250, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(250, 3).Value = "Metropolitana"
$ws.Cells.Item(250, 4).Value = 44627
$ws.Cells.Item(250, 5).Value = 13
$ws.Cells.Item(250, 6).Value = "Fruta"
$ws.Cells.Item(250, 7).Value = 100101
$ws.Cells.Item(250, 8).Value = "Berries"
$ws.Cells.Item(250, 9).Value = 100101001
$ws.Cells.Item(250, 10).Value = "Arándano (blue)"
$ws.Cells.Item(250, 11).Value = "Sin especificar"
$ws.Cells.Item(250, 12).Value = "Primera"
$ws.Cells.Item(250, 13).Value = 1500
$ws.Cells.Item(250, 14).Value = 4000
$ws.Cells.Item(250, 15).Value = 4000
$ws.Cells.Item(250, 16).Value = 4000
$ws.Cells.Item(250, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(250, 18).Value = "Provincia de Linares"
$ws.Cells.Item(250, 19).Value = 2000
$ws.Cells.Item(250, 20).Value = 2
